{"js": "// Office.js (Word JavaScript API) edit script\n// Applies three changes described by the diff:\n// 1) \"21 years\" -> \"15+ years\" in the Professional Summary paragraph\n// 2) Rewrite the FLEEM bullet under \"RESEARCH DIRECTOR - Progressive Change Campaign Committee\"\n// 3) Insert a new bullet after \"Developed innovative approaches to visualizing demographic...\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst OLD_SUMMARY_FRAGMENT = \"21 years of experience\";\nconst NEW_SUMMARY_FRAGMENT = \"15+ years of experience\";\n\nconst OLD_FLEEM =\n  \"\\u2022 Engineered FLEEM web application using Twilio's API to make thousands of simultaneous phone calls for IVR polls\";\nconst NEW_FLEEM =\n  \"\\u2022 Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys\";\n\nconst ANCHOR_TEXT =\n  \"\\u2022 Developed innovative approaches to visualizing demographic and market data, enhancing clients' understanding of research findings\";\nconst NEW_BULLET_TEXT =\n  \"\\u2022 Trained staff on building Python tooling for report generation and analysis\";\n\nlet anchorParagraph = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n\n  if (text.indexOf(OLD_SUMMARY_FRAGMENT) !== -1) {\n    // Replace just the \"21 years\" -> \"15+ years\" fragment, preserving the rest of the text.\n    const range = para.getRange();\n    const results = range.search(OLD_SUMMARY_FRAGMENT, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n    if (results.items.length > 0) {\n      results.items[0].insertText(NEW_SUMMARY_FRAGMENT, Word.InsertLocation.replace);\n    }\n  } else if (text === OLD_FLEEM) {\n    para.insertText(NEW_FLEEM, Word.InsertLocation.replace);\n  } else if (text === ANCHOR_TEXT) {\n    anchorParagraph = para;\n  }\n}\n\nawait context.sync();\n\nif (anchorParagraph) {\n  anchorParagraph.insertParagraph(NEW_BULLET_TEXT, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script\n# Applies three changes described by the diff:\n# 1) \"21 years\" -> \"15+ years\" in the Professional Summary paragraph\n# 2) Rewrite the FLEEM bullet under \"RESEARCH DIRECTOR - Progressive Change Campaign Committee\"\n# 3) Insert a new bullet after \"Developed innovative approaches to visualizing demographic...\"\n\n$d = $word.ActiveDocument\n$bullet = [char]0x2022\n\n# --- 1) Update experience years in the Professional Summary ---\n$find1 = $d.Content.Find\n$find1.Execute(\n    \"21 years of experience\",\n    $false, $true, $false, $false, $false, $true, 1, $false,\n    \"15+ years of experience\", 2\n)\n\n# --- 2) Rewrite the FLEEM bullet (Research Director section) ---\n$find2 = $d.Content.Find\n$oldFleem = \"Engineered FLEEM web application using Twilio's API to make thousands of simultaneous phone calls for IVR polls\"\n$newFleem = \"Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys\"\n$find2.Execute(\n    $oldFleem,\n    $false, $true, $false, $false, $false, $true, 1, $false,\n    $newFleem, 2\n)\n\n# --- 3) Insert new bullet after the \"visualizing demographic\" bullet ---\n$anchorText = \"Developed innovative approaches to visualizing demographic and market data, enhancing clients' understanding of research findings\"\n$newBulletText = \"$bullet Trained staff on building Python tooling for report generation and analysis\"\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*$anchorText*\") {\n        $p.Range.InsertParagraphAfter()\n        $newPara = $p.Next()\n        $newPara.Range.Text = $newBulletText\n        break\n    }\n}\n"}
